$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.056.99"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "3.587.06"
$ws.Range("E3").Value = "  -2.61%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "2.28"
$ws.Range("E5").Value = "  +18.26%  "

$ws.Range("D6").Value = "225.02"
$ws.Range("E6").Value = "  -4.76%  "

$ws.Range("D7").Value = "634.29"
$ws.Range("E7").Value = "  -3.03%  "

$ws.Range("D8").Value = "0.408"
$ws.Range("E8").Value = "  -3.32%  "

$ws.Range("D9").Value = "1.08"
$ws.Range("E9").Value = "  +0.84%  "

$ws.Range("E10").Value = "  +0.16%  "

$ws.Range("D11").Value = "3.584.91"
$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("D12").Value = "45.92"
$ws.Range("E12").Value = "  +4.36%  "

$ws.Range("D13").Value = "0.205"
$ws.Range("E13").Value = "  -1.54%  "

$ws.Range("D14").Value = "0.0000287"
$ws.Range("E14").Value = "  -3.37%  "

$ws.Range("D15").Value = "6.43"
$ws.Range("E15").Value = "  -4.49%  "

$ws.Range("D16").Value = "4.267.32"
$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").Value = "94.838.67"
$ws.Range("E17").Value = "  -1.73%  "

$ws.Range("D18").Value = "8.73"
$ws.Range("E18").Value = "  -2.61%  "

$ws.Range("D19").Value = "20.21"
$ws.Range("E19").Value = "  +8.85%  "

$ws.Range("D20").Value = "3.598.32"
$ws.Range("E20").Value = "  -2.31%  "

$ws.Range("D21").Value = "12.84"
$ws.Range("E21").Value = "  -0.94%  "

$ws.Range("D22").Value = "0.507"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").Value = "499.22"
$ws.Range("E23").Value = "  -4.06%  "

$ws.Range("D24").Value = "3.21"
$ws.Range("E24").Value = "  -5.49%  "

$ws.Range("D25").Value = "0.239"
$ws.Range("E25").Value = "  +17.87%  "

$ws.Range("D26").Value = "116.23"
$ws.Range("E26").Value = "  +14.91%  "

$ws.Range("D27").Value = "0.0000201"
$ws.Range("E27").Value = "  -4.30%  "

$ws.Range("D28").Value = "6.71"
$ws.Range("E28").Value = "  -2.47%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.791.54"
$ws.Range("E29").Value = "  -2.41%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "12.50"
$ws.Range("E30").Value = "  -6.16%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "12.67"
$ws.Range("E31").Value = "  +2.75%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.88"
$ws.Range("E32").Value = "  -3.77%  "

$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").Value = "  -5.27%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "1.76"
$ws.Range("E36").Value = "  -4.27%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "31.49"
$ws.Range("E37").Value = "  -1.82%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.580"
$ws.Range("E38").Value = "  -1.48%  "

$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "587.25"
$ws.Range("E40").Value = "  -8.52%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "8.24"
$ws.Range("E41").Value = "  -5.79%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "6.73"
$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "40.69"
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.157"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.468"
$ws.Range("E45").Value = "  +2.89%  "

$ws.Range("D46").Value = "1.90"
$ws.Range("E46").Value = "  -6.71%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0467"
$ws.Range("E47").Value = "  +2.83%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "0.914"
$ws.Range("E48").Value = "  -3.72%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "23.43"
$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "8.48"
$ws.Range("E50").Value = "  -0.12%  "

$ws.Range("B51").Value = "MantraDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D51").Value = "3.58"
$ws.Range("E51").Value = "  +1.44%  "
